$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 34486050
$ws.Range("I98").Value = 37039908
$ws.Range("J98").Value = 8953
$ws.Range("K98").Value = 37039908
$ws.Range("L98").Value = 8953
$ws.Range("M98").Value = -37038410
$ws.Range("N98").Value = -11949
$ws.Range("H122").Value = 34486050
$ws.Range("I122").Value = 37039908
$ws.Range("J122").Value = 8953
$ws.Range("K122").Value = 111119724
$ws.Range("L122").Value = 26859
$ws.Range("M122").Value = -111117274
$ws.Range("N122").Value = -31759
$ws.Range("H123").Value = 70750.5
$ws.Range("J123").Value = 70750.5
$ws.Range("L123").Value = 70750.5
$ws.Range("N123").Value = -80550.5
$ws.Range("H137").Value = 2948.9678
$ws.Range("I137").Value = 2894.5454
$ws.Range("K137").Value = 8683.636200000001
$ws.Range("M137").Value = -6133.636200000001
$ws.Range("H138").Value = 4778.896
$ws.Range("J138").Value = 8938.781999999999
$ws.Range("L138").Value = 26816.346
$ws.Range("N138").Value = -37096.346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3910867
$ws.Range("I32").Value = 4314785.5
$ws.Range("K32").Value = 4314785.5
$ws.Range("M32").Value = -4314498.5
$ws.Range("H42").Value = 19663.334
$ws.Range("I42").Value = 19500
$ws.Range("J42").Value = 19990
$ws.Range("K42").Value = 19500
$ws.Range("L42").Value = 19990
$ws.Range("M42").Value = -19014
$ws.Range("N42").Value = -20962
$ws.Range("H61").Value = 3933.0476
$ws.Range("I61").Value = 2893.8125
$ws.Range("K61").Value = 2893.8125
$ws.Range("M61").Value = -2681.8125
$ws.Range("H63").Value = 1828.3334
$ws.Range("J63").Value = 1860
$ws.Range("L63").Value = 1860
$ws.Range("N63").Value = -3232
$ws.Range("H66").Value = 1828.3334
$ws.Range("J66").Value = 1860
$ws.Range("L66").Value = 9300
$ws.Range("N66").Value = -16164
$ws.Range("H74").Value = 42957.562
$ws.Range("I74").Value = 51253.656
$ws.Range("K74").Value = 51253.656
$ws.Range("M74").Value = -50379.656
$ws.Range("H77").Value = 42957.562
$ws.Range("I77").Value = 51253.656
$ws.Range("K77").Value = 256268.28
$ws.Range("M77").Value = -251900.28
$ws.Range("H122").Value = 9697.333000000001
$ws.Range("I122").Value = 10236.32
$ws.Range("K122").Value = 30708.96
$ws.Range("M122").Value = -28258.96
$ws.Range("H132").Value = 7607.8423
$ws.Range("I132").Value = 6125.4165
$ws.Range("K132").Value = 18376.2495
$ws.Range("M132").Value = -15846.2495
$ws.Range("H136").Value = 3933.0476
$ws.Range("I136").Value = 2893.8125
$ws.Range("K136").Value = 8681.4375
$ws.Range("M136").Value = -6131.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2873.7354
$ws.Range("I105").Value = 2368.08
$ws.Range("J105").Value = 4278.3335
$ws.Range("K105").Value = 2368.08
$ws.Range("L105").Value = 4278.3335
$ws.Range("M105").Value = -621.0799999999999
$ws.Range("N105").Value = -7772.3335
$ws.Range("H134").Value = 3697.2856
$ws.Range("I134").Value = 1952.3175
$ws.Range("K134").Value = 5856.9525
$ws.Range("M134").Value = -3321.9525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = $null
$ws.Range("N6").Value = 0
$ws.Range("H31").Value = 5996.923
$ws.Range("I31").Value = 2594.0652
$ws.Range("J31").Value = 10888.531
$ws.Range("K31").Value = 2594.0652
$ws.Range("L31").Value = 10888.531
$ws.Range("M31").Value = -2299.0652
$ws.Range("N31").Value = -11478.531
$ws.Range("H34").Value = 5996.923
$ws.Range("I34").Value = 2594.0652
$ws.Range("J34").Value = 10888.531
$ws.Range("K34").Value = 2594.0652
$ws.Range("L34").Value = 10888.531
$ws.Range("M34").Value = -2392.0652
$ws.Range("N34").Value = -11292.531
$ws.Range("H62").Value = 5995
$ws.Range("I62").Value = 4993.75
$ws.Range("K62").Value = 4993.75
$ws.Range("M62").Value = -4369.75
$ws.Range("H65").Value = 5995
$ws.Range("I65").Value = 4993.75
$ws.Range("K65").Value = 24968.75
$ws.Range("M65").Value = -21848.75
$ws.Range("H99").Value = 11560.8
$ws.Range("I99").Value = 27450
$ws.Range("J99").Value = 7588.5
$ws.Range("K99").Value = 27450
$ws.Range("L99").Value = 7588.5
$ws.Range("M99").Value = -25952
$ws.Range("N99").Value = -10584.5
$ws.Range("H118").Value = 94913.664
$ws.Range("J118").Value = 94913.664
$ws.Range("L118").Value = 94913.664
$ws.Range("N118").Value = -98227.664
$ws.Range("H119").Value = 96000
$ws.Range("J119").Value = 96000
$ws.Range("L119").Value = 96000
$ws.Range("N119").Value = -105676
$ws.Range("H122").Value = 1385.2307
$ws.Range("I122").Value = 1141.5555
$ws.Range("K122").Value = 3424.6665
$ws.Range("M122").Value = -974.6664999999998
$ws.Range("H126").Value = 11560.8
$ws.Range("I126").Value = 27450
$ws.Range("J126").Value = 7588.5
$ws.Range("K126").Value = 82350
$ws.Range("L126").Value = 22765.5
$ws.Range("M126").Value = -79880
$ws.Range("N126").Value = -27705.5
$ws.Range("H132").Value = 5735.524
$ws.Range("I132").Value = 1870.7778
$ws.Range("K132").Value = 5612.3334
$ws.Range("M132").Value = -3082.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3185.8235
$ws.Range("J113").Value = 3345.3333
$ws.Range("L113").Value = 10035.9999
$ws.Range("N113").Value = -14375.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 43384.285
$ws.Range("J58").Value = 73724.75
$ws.Range("L58").Value = 73724.75
$ws.Range("N58").Value = -74278.75
$ws.Range("H113").Value = 8069
$ws.Range("I113").Value = 4241.5
$ws.Range("J113").Value = 9600
$ws.Range("K113").Value = 4241.5
$ws.Range("L113").Value = 9600
$ws.Range("M113").Value = -2071.5
$ws.Range("N113").Value = -13940
$ws.Range("H122").Value = 41315.85
$ws.Range("I122").Value = 80165.766
$ws.Range("J122").Value = 5240.9287
$ws.Range("K122").Value = 240497.298
$ws.Range("L122").Value = 15722.7861
$ws.Range("M122").Value = -238047.298
$ws.Range("N122").Value = -20622.7861
$ws.Range("H126").Value = 2655.125
$ws.Range("I126").Value = 2412.7144
$ws.Range("K126").Value = 7238.1432
$ws.Range("M126").Value = -4768.1432
$ws.Range("H132").Value = 2363.342
$ws.Range("I132").Value = 2314.5
$ws.Range("K132").Value = 6943.5
$ws.Range("M132").Value = -4413.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5866.85
$ws.Range("I40").Value = 5564.154
$ws.Range("J40").Value = 6429
$ws.Range("K40").Value = 5564.154
$ws.Range("L40").Value = 6429
$ws.Range("M40").Value = -5428.154
$ws.Range("N40").Value = -6701
$ws.Range("H46").Value = 16178407
$ws.Range("I46").Value = 11494718
$ws.Range("J46").Value = 18520252
$ws.Range("K46").Value = 11494718
$ws.Range("L46").Value = 18520252
$ws.Range("M46").Value = -11494530
$ws.Range("N46").Value = -18520628
$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488
$ws.Range("H82").Value = 1836.1052
$ws.Range("I82").Value = 1599.2
$ws.Range("J82").Value = 2099.3333
$ws.Range("K82").Value = 1599.2
$ws.Range("L82").Value = 2099.3333
$ws.Range("M82").Value = -1238.2
$ws.Range("N82").Value = -2821.3333
$ws.Range("H85").Value = 1836.1052
$ws.Range("I85").Value = 1599.2
$ws.Range("J85").Value = 2099.3333
$ws.Range("K85").Value = 1599.2
$ws.Range("L85").Value = 2099.3333
$ws.Range("M85").Value = -351.2
$ws.Range("N85").Value = -4595.3333
$ws.Range("H132").Value = 17864404
$ws.Range("I132").Value = 33337864
$ws.Range("J132").Value = 10410.615
$ws.Range("K132").Value = 100013592
$ws.Range("L132").Value = 31231.845
$ws.Range("M132").Value = -100011062
$ws.Range("N132").Value = -36291.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 15416.667
$ws.Range("J43").Value = 15416.667
$ws.Range("L43").Value = 15416.667
$ws.Range("N43").Value = -15714.667
$ws.Range("H81").Value = 28589928
$ws.Range("I81").Value = 1833
$ws.Range("K81").Value = 3666
$ws.Range("M81").Value = -2605
$ws.Range("H84").Value = 28589928
$ws.Range("I84").Value = 1833
$ws.Range("K84").Value = 18330
$ws.Range("M84").Value = -13026
$ws.Range("H132").Value = 71517064
$ws.Range("I132").Value = 125028750
$ws.Range("J132").Value = 168148
$ws.Range("K132").Value = 375086250
$ws.Range("L132").Value = 504444
$ws.Range("M132").Value = -375083720
$ws.Range("N132").Value = -509504
$ws.Range("H136").Value = 34486010
$ws.Range("I136").Value = 52632310
$ws.Range("K136").Value = 157896930
$ws.Range("M136").Value = -157894380

Write-Output "Applied all cell updates"